$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 76, pushing existing data (old rows 76-129) down to 78-131.
$ws.Rows.Item(76).Insert()
$ws.Rows.Item(76).Insert()

# Populate new row 76 (Modesto / Primera)
$ws.Range("A76").Value2 = 8
$ws.Range("B76").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C76").Value2 = "Coquimbo"
$ws.Range("D76").Value2 = 44957
$ws.Range("E76").Value2 = 4
$ws.Range("F76").Value2 = "Fruta"
$ws.Range("G76").Value2 = 100103
$ws.Range("H76").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I76").Value2 = 100103003
$ws.Range("J76").Value2 = "Damasco"
$ws.Range("K76").Value2 = "Modesto"
$ws.Range("L76").Value2 = "Primera"
$ws.Range("M76").Value2 = 200
$ws.Range("N76").Value2 = 21000
$ws.Range("O76").Value2 = 22000
$ws.Range("P76").Value2 = 21500
$ws.Range("Q76").Value2 = "$/caja 16 kilos"
$ws.Range("R76").Value2 = "Región Metropolitana"
$ws.Range("S76").Value2 = 1344
$ws.Range("T76").Value2 = 16

# Populate new row 77 (Modesto / Segunda)
$ws.Range("A77").Value2 = 8
$ws.Range("B77").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value2 = "Coquimbo"
$ws.Range("D77").Value2 = 44957
$ws.Range("E77").Value2 = 4
$ws.Range("F77").Value2 = "Fruta"
$ws.Range("G77").Value2 = 100103
$ws.Range("H77").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I77").Value2 = 100103003
$ws.Range("J77").Value2 = "Damasco"
$ws.Range("K77").Value2 = "Modesto"
$ws.Range("L77").Value2 = "Segunda"
$ws.Range("M77").Value2 = 160
$ws.Range("N77").Value2 = 17000
$ws.Range("O77").Value2 = 18000
$ws.Range("P77").Value2 = 17500
$ws.Range("Q77").Value2 = "$/caja 16 kilos"
$ws.Range("R77").Value2 = "Región Metropolitana"
$ws.Range("S77").Value2 = 1094
$ws.Range("T77").Value2 = 16
